# Rename variables for consistency with standard:
#   - "soil_descriptors" sheet: header A1 "class"     -> "soil_class"
#   - "soil_spectra"     sheet: header A1 "soil_type"  -> "soil_class"

$wb = $excel.ActiveWorkbook

$wsDescriptors = $wb.Worksheets.Item("soil_descriptors")
$wsSpectra = $wb.Worksheets.Item("soil_spectra")

$wsDescriptors.Range("A1").Value = "soil_class"
$wsSpectra.Range("A1").Value = "soil_class"

# Reproduce the final selection/active-sheet state recorded in the workbook:
# the user last worked on soil_spectra (cell F2 selected) then switched back
# to soil_descriptors (cell A2 selected), leaving soil_descriptors active.
$wsSpectra.Activate()
$wsSpectra.Range("F2").Select()

$wsDescriptors.Activate()
$wsDescriptors.Range("A2").Select()
